$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 25
$ws.Range("H4").Value = 25

$ws.Range("E5").Value = 117
$ws.Range("F5").Value = 73
$ws.Range("H5").Value = 73

$ws.Range("F6").Value = 25
$ws.Range("H6").Value = 25

$ws.Range("E10").Value = 390
$ws.Range("F10").Value = 187
$ws.Range("H10").Value = 187

$ws.Range("F11").Value = 139
$ws.Range("H11").Value = 139

$ws.Range("E12").Value = 375
$ws.Range("F12").Value = 208
$ws.Range("H12").Value = 208

$ws.Range("E14").Value = 98
$ws.Range("F14").Value = 49
$ws.Range("H14").Value = 49

$ws.Range("E15").Value = 124
$ws.Range("F15").Value = 49
$ws.Range("H15").Value = 49

$ws.Range("F16").Value = 81
$ws.Range("H16").Value = 81

$ws.Range("E17").Value = 73
$ws.Range("F17").Value = 35
$ws.Range("H17").Value = 35

$ws.Range("F21").Value = 65
$ws.Range("H21").Value = 65

$ws.Range("F22").Value = 71
$ws.Range("H22").Value = 71

$ws.Range("E23").Value = 167
$ws.Range("F23").Value = 73
$ws.Range("H23").Value = 73

$ws.Range("F24").Value = 85
$ws.Range("H24").Value = 85

$ws.Range("E25").Value = 194
$ws.Range("F25").Value = 89
$ws.Range("H25").Value = 89

$ws.Range("E26").Value = 114
$ws.Range("F26").Value = 69
$ws.Range("H26").Value = 69

$ws.Range("E27").Value = 257
$ws.Range("F27").Value = 123
$ws.Range("H27").Value = 123

$ws.Range("F28").Value = 50
$ws.Range("H28").Value = 50

$ws.Range("E29").Value = 141
$ws.Range("F29").Value = 79
$ws.Range("H29").Value = 79

$ws.Range("F31").Value = 29
$ws.Range("H31").Value = 29

$ws.Range("F32").Value = 88
$ws.Range("H32").Value = 88

$ws.Range("E33").Value = 237
$ws.Range("F33").Value = 121
$ws.Range("H33").Value = 121

$ws.Range("F34").Value = 102
$ws.Range("H34").Value = 102

$ws.Range("F35").Value = 72
$ws.Range("H35").Value = 72

$ws.Range("F36").Value = 31
$ws.Range("H36").Value = 31

$ws.Range("E37").Value = 126
$ws.Range("F37").Value = 61
$ws.Range("H37").Value = 61

$ws.Range("F38").Value = 51
$ws.Range("H38").Value = 51

$ws.Range("E39").Value = 158
$ws.Range("F39").Value = 75
$ws.Range("H39").Value = 75

$ws.Range("E40").Value = 214
$ws.Range("F40").Value = 95
$ws.Range("H40").Value = 95

$ws.Range("F41").Value = 137
$ws.Range("H41").Value = 137

$ws.Range("F42").Value = 151
$ws.Range("H42").Value = 151

$ws.Range("E43").Value = 95
$ws.Range("F43").Value = 48
$ws.Range("H43").Value = 48

$ws.Range("F44").Value = 119
$ws.Range("H44").Value = 119

$ws.Range("F45").Value = 50
$ws.Range("H45").Value = 50

$ws.Range("F46").Value = 136
$ws.Range("H46").Value = 136

$ws.Range("F47").Value = 172
$ws.Range("H47").Value = 172

$ws.Range("E48").Value = 166
$ws.Range("F48").Value = 66
$ws.Range("H48").Value = 66

$ws.Range("E49").Value = 236
$ws.Range("F49").Value = 101
$ws.Range("H49").Value = 101

$ws.Range("E50").Value = 207
$ws.Range("F50").Value = 82
$ws.Range("H50").Value = 82

$ws.Range("E51").Value = 191
$ws.Range("F51").Value = 79
$ws.Range("H51").Value = 79

$ws.Range("F52").Value = 11
$ws.Range("H52").Value = 11

